# Auto-generated: apply scheduled-runner price/profit refresh to Sheets workbook
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 29665.666
$ws.Range("I44").Value = 19998
$ws.Range("J44").Value = 34499.5
$ws.Range("K44").Value = 19998
$ws.Range("L44").Value = 34499.5
$ws.Range("M44").Value = -19536
$ws.Range("N44").Value = -35423.5
$ws.Range("H51").Value = 4067.2666
$ws.Range("I51").Value = 1866.6666
$ws.Range("J51").Value = 4617.4165
$ws.Range("K51").Value = 1866.6666
$ws.Range("L51").Value = 4617.4165
$ws.Range("M51").Value = -1382.6666
$ws.Range("N51").Value = -5585.4165
$ws.Range("H62").Value = 7678.091
$ws.Range("I62").Value = 7773.8887
$ws.Range("K62").Value = 7773.8887
$ws.Range("M62").Value = -7149.8887
$ws.Range("H65").Value = 7678.091
$ws.Range("I65").Value = 7773.8887
$ws.Range("K65").Value = 38869.4435
$ws.Range("M65").Value = -35749.4435
$ws.Range("H87").Value = 200000
$ws.Range("J87").Value = 200000
$ws.Range("L87").Value = 200000
$ws.Range("N87").Value = -202496
$ws.Range("H90").Value = 200000
$ws.Range("J90").Value = 200000
$ws.Range("L90").Value = 600000
$ws.Range("N90").Value = -612480
$ws.Range("H111").Value = 21580.5
$ws.Range("I111").Value = 1029
$ws.Range("J111").Value = 42132
$ws.Range("K111").Value = 3087
$ws.Range("L111").Value = 126396
$ws.Range("M111").Value = -20
$ws.Range("N111").Value = -132530
$ws.Range("H129").Value = 1268.0834
$ws.Range("J129").Value = 2575.7144
$ws.Range("L129").Value = 7727.1432
$ws.Range("N129").Value = -17727.1432
$ws.Range("H132").Value = 3636.0322
$ws.Range("I132").Value = 3311.1072
$ws.Range("J132").Value = 6668.6665
$ws.Range("K132").Value = 9933.321599999999
$ws.Range("L132").Value = 20005.9995
$ws.Range("M132").Value = -7403.321599999999
$ws.Range("N132").Value = -25065.9995
$ws.Range("H137").Value = 20841666
$ws.Range("I137").Value = 25009316
$ws.Range("K137").Value = 75027948
$ws.Range("M137").Value = -75025398

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11601.64
$ws.Range("I32").Value = 6774.22
$ws.Range("K32").Value = 6774.22
$ws.Range("M32").Value = -6487.22
$ws.Range("H45").Value = 1541.2333
$ws.Range("I45").Value = 1530.1333
$ws.Range("K45").Value = 1530.1333
$ws.Range("M45").Value = -1153.1333

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4380
$ws.Range("I20").Value = 3831.3333
$ws.Range("K20").Value = 3831.3333
$ws.Range("M20").Value = -3584.3333
$ws.Range("H134").Value = 1563.5217
$ws.Range("I134").Value = 1045.762
$ws.Range("K134").Value = 3137.286
$ws.Range("M134").Value = -602.2860000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1783.0769
$ws.Range("I16").Value = 1515.5
$ws.Range("K16").Value = 1515.5
$ws.Range("M16").Value = -1228.5
$ws.Range("H31").Value = 376686.03
$ws.Range("I31").Value = 10770.385
$ws.Range("J31").Value = 566962.2
$ws.Range("K31").Value = 10770.385
$ws.Range("L31").Value = 566962.2
$ws.Range("M31").Value = -10475.385
$ws.Range("N31").Value = -567552.2
$ws.Range("H34").Value = 376686.03
$ws.Range("I34").Value = 10770.385
$ws.Range("J34").Value = 566962.2
$ws.Range("K34").Value = 10770.385
$ws.Range("L34").Value = 566962.2
$ws.Range("M34").Value = -10568.385
$ws.Range("N34").Value = -567366.2
$ws.Range("H113").Value = 1783.0769
$ws.Range("I113").Value = 1515.5
$ws.Range("K113").Value = 1515.5
$ws.Range("M113").Value = 654.5
$ws.Range("H122").Value = 2837.4443
$ws.Range("J122").Value = 2991.75
$ws.Range("L122").Value = 8975.25
$ws.Range("N122").Value = -13875.25
$ws.Range("H132").Value = 1860.2632
$ws.Range("I132").Value = 1769.1945
$ws.Range("K132").Value = 5307.583500000001
$ws.Range("M132").Value = -2777.583500000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1350
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H127").Value = 802
$ws.Range("J127").Value = 802
$ws.Range("L127").Value = 2406
$ws.Range("N127").Value = -12326
$ws.Range("H131").Value = 5835.5293
$ws.Range("I131").Value = 2899.5
$ws.Range("J131").Value = 6227
$ws.Range("K131").Value = 8698.5
$ws.Range("L131").Value = 18681
$ws.Range("M131").Value = -3658.5
$ws.Range("N131").Value = -28761
$ws.Range("H136").Value = 3299.5
$ws.Range("I136").Value = 3299.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9898.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4798.5
$ws.Range("N136").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3031.7812
$ws.Range("I102").Value = 2214.125
$ws.Range("K102").Value = 2214.125
$ws.Range("M102").Value = -592.125
$ws.Range("H113").Value = 3081.0833
$ws.Range("I113").Value = 2800.6
$ws.Range("J113").Value = 3281.4285
$ws.Range("K113").Value = 2800.6
$ws.Range("L113").Value = 3281.4285
$ws.Range("M113").Value = -630.5999999999999
$ws.Range("N113").Value = -7621.4285
$ws.Range("H122").Value = 4142.8237
$ws.Range("I122").Value = 3887.8572
$ws.Range("K122").Value = 11663.5716
$ws.Range("M122").Value = -9213.571599999999
$ws.Range("H132").Value = 4631.129
$ws.Range("I132").Value = 4621.6
$ws.Range("K132").Value = 13864.8
$ws.Range("M132").Value = -11334.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10417747
$ws.Range("I46").Value = 20834494
$ws.Range("K46").Value = 20834494
$ws.Range("M46").Value = -20834306
$ws.Range("H122").Value = 8350.75
$ws.Range("I122").Value = 5614.3335
$ws.Range("K122").Value = 16843.0005
$ws.Range("M122").Value = -14393.0005
$ws.Range("H136").Value = 5625.8
$ws.Range("I136").Value = 6836.6875
$ws.Range("K136").Value = 20510.0625
$ws.Range("M136").Value = -17960.0625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 19794610
$ws.Range("I122").Value = 17048266
$ws.Range("J122").Value = 50004396
$ws.Range("K122").Value = 51144798
$ws.Range("L122").Value = 150013188
$ws.Range("M122").Value = -51142348
$ws.Range("N122").Value = -150018088
$ws.Range("H132").Value = 3591.3809
$ws.Range("I132").Value = 2577.0908
$ws.Range("K132").Value = 7731.2724
$ws.Range("M132").Value = -5201.2724

